# IO Updates to model, updated with current US .mdl
#
# 1) "SoCaOMSbRIC" sheet: the combined "ISIC 20T21" column is split into two
#    separate columns, "ISIC 20" and "ISIC 21".
# 2) "About" sheet: the Notes paragraph explaining capital/OM spending
#    allocation is reworded, and the process-emissions sentence is replaced
#    with a note about a different input variable.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# SoCaOMSbRIC sheet - split the "ISIC 20T21" column into "ISIC 20" / "ISIC 21"
# ---------------------------------------------------------------------
$wsMain = $wb.Worksheets.Item("SoCaOMSbRIC")

# Insert a new column before the existing "ISIC 20T21" column (column K).
$wsMain.Range("K1").EntireColumn.Insert()

# Relabel the two now-adjacent header cells.
$wsMain.Range("K1").Value = "ISIC 20"
$wsMain.Range("L1").Value = "ISIC 21"

# Both new columns carry the same (zero) share value the combined column had.
$wsMain.Range("K2").Value = 0
$wsMain.Range("L2").Value = 0

# ---------------------------------------------------------------------
# About sheet - update the Notes section text
# ---------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A21").Value = "to cover the allocation of capital and OM expensies for industry energy efficiency and CCS."
$wsAbout.Range("A22").ClearContents()
$wsAbout.Range("A23").Value = "A different input variable governs the breakdown of revenues due to Process Emissions policies."
$wsAbout.Range("A24").ClearContents()
